# Auto-generated Excel COM-interop script
# Applies the "nominal nas filling repaired" fix across the 3 sheets:
#   numeric_variable, object_variable, column with soo000oo000oo000oo000oo long name

$wb = $excel.ActiveWorkbook

# ---- numeric_variable ----
$ws = $wb.Worksheets.Item("numeric_variable")
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2
$ws.Range("N4").Value = 0.6783525101020478
$ws.Range("P4").Value = 0.3567050202040956
$ws.Range("Q4").Value = 0.5754920420347929
$ws.Range("S4").Value = 0.1509840840695857
$ws.Range("T4").Value = 0.7796336996336996
$ws.Range("V4").Value = 0.5592673992673991
$ws.Range("B15").Value = 9.85252143345858
$ws.Range("B16").Value = 5.020494655354245
$ws.Range("B17").Value = -9.351980052667205
$ws.Range("B18").Value = 6.466083189783197
$ws.Range("B19").Value = 10.07544999943667
$ws.Range("B20").Value = 13.38880928291256
$ws.Range("B21").Value = 29.12869353257992

# ---- object_variable ----
$ws = $wb.Worksheets.Item("object_variable")
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2
$ws.Range("N4").Value = 0.9502281521813574
$ws.Range("P4").Value = 0.9004563043627147
$ws.Range("Q4").Value = 0.8755411255411255
$ws.Range("S4").Value = 0.751082251082251
$ws.Range("T4").Value = 0.9302380952380953
$ws.Range("V4").Value = 0.8604761904761906
# Swap the "2.0%"/"1.0%" counts and their display text
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("F13").Value = "'1.0%"
$ws.Range("H13").Value = "'2.0%"

# Rewrite the value_counts table (rows 14-24)
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 113
$ws.Range("C14").Value = 106
$ws.Range("D14").Value = 93.80530973451327
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 6.194690265486726
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = 110
$ws.Range("C15").Value = 19
$ws.Range("D15").Value = 17.27272727272727
$ws.Range("E15").Value = 90
$ws.Range("F15").Value = 81.81818181818181
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.9090909090909091
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = 104
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5.769230769230769
$ws.Range("E16").Value = 98
$ws.Range("F16").Value = 94.23076923076923
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = 103
$ws.Range("C17").Value = 87
$ws.Range("D17").Value = 84.46601941747574
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 15.53398058252427
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = 102
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 77
$ws.Range("F18").Value = 75.49019607843137
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = 24.50980392156863
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = 95
$ws.Range("C19").Value = 60
$ws.Range("D19").Value = 63.1578947368421
$ws.Range("E19").Value = 35
$ws.Range("F19").Value = 36.8421052631579
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 92
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 87
$ws.Range("F20").Value = 94.56521739130434
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 5.434782608695652
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = 91
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = 65.93406593406593
$ws.Range("G21").Value = 31
$ws.Range("H21").Value = 34.06593406593407
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = 89
$ws.Range("C22").Value = 36
$ws.Range("D22").Value = 40.44943820224719
$ws.Range("E22").Value = 53
$ws.Range("F22").Value = 59.55056179775281
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("A23").Value = 0
$ws.Range("B23").Value = 58
$ws.Range("C23").Value = 58
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = 43
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 34.88372093023256
$ws.Range("G24").Value = 28
$ws.Range("H24").Value = 65.11627906976744

# ---- column with soo000oo000oo000oo0 ----
$ws = $wb.Worksheets.Item("column with soo000oo000oo000oo0")
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2
$ws.Range("N4").Value = 0.5022472775837271
$ws.Range("P4").Value = 0.004494555167454228
$ws.Range("Q4").Value = 0.5062682051529634
$ws.Range("S4").Value = 0.01253641030592689
$ws.Range("T4").Value = 0.5126129426129427
$ws.Range("V4").Value = 0.02522588522588531
$ws.Range("B15").Value = 10.00404639041234
$ws.Range("B16").Value = 4.780183502172311
$ws.Range("B17").Value = -4.678935533229469
$ws.Range("B18").Value = 6.677748398500486
$ws.Range("B19").Value = 10.11121159172304
$ws.Range("B20").Value = 13.12318206170177
$ws.Range("B21").Value = 24.14265030523879
